$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row is the next one after the current last data row (row 73 -> new row 74)
$newRow = 74
$prevRow = $newRow - 1

# Copy formatting from the previous row so the new row picks up the same
# cell styles (bold/border index style on column A, date style on column E).
$ws.Range("A" + $prevRow + ":V" + $prevRow).Copy() | Out-Null
$ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 73
$ws.Cells.Item($newRow, 2).Value = "armenia"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45237.625
$ws.Cells.Item($newRow, 6).Value = "Noah"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Ararat-Armenia"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 3.14
$ws.Cells.Item($newRow, 11).Value = "06/11/2023 03:12"
$ws.Cells.Item($newRow, 12).Value = 3.03
$ws.Cells.Item($newRow, 13).Value = "07/11/2023 14:55"
$ws.Cells.Item($newRow, 14).Value = 3.36
$ws.Cells.Item($newRow, 15).Value = "06/11/2023 03:12"
$ws.Cells.Item($newRow, 16).Value = 3.69
$ws.Cells.Item($newRow, 17).Value = "07/11/2023 14:58"
$ws.Cells.Item($newRow, 18).Value = 2.06
$ws.Cells.Item($newRow, 19).Value = "06/11/2023 03:12"
$ws.Cells.Item($newRow, 20).Value = 2.19
$ws.Cells.Item($newRow, 21).Value = "07/11/2023 14:58"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/noah-ararat-armenia/ADdhtgKh/"
